$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 13.1321
$ws.Range("E4").Value = 13.55529999999998
$ws.Range("E7").Value = 11.8696
$ws.Range("E8").Value = 13.35809999999999
$ws.Range("A11").Value = -21.90820000000003
$ws.Range("A12").Value = -22.86990000000002
$ws.Range("E12").Value = 12.53679999999999
$ws.Range("E14").Value = 13.63080000000001
$ws.Range("A15").Value = -21.51350000000003
$ws.Range("E22").Value = 12.4822
